$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds date serial values that need to move forward
# by one day, from 45207 (2023-10-08) to 45208 (2023-10-09), for every data
# row currently set to 45207.
for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value = 45208
    }
}
